$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "Relates To ID" value for row 3 (VIS_002 segment)
$ws.Range("J3").Value = "BE_001, BE_002"

# Update the active selection to reflect where the user ended up (J7)
$ws.Range("J7").Select()
